$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.036099460753967
$bf[0,2] = 1.042337898687244
$bf[0,3] = 1.044645817576086
$bf[0,4] = 1.05455329833542
$ws.Range("B2:F2").Value = $bf

$inn = New-Object 'object[,]' 1,6
$inn[0,0] = 1.034184161016004
$inn[0,1] = 1.041209624849915
$inn[0,2] = 1.04511486930196
$inn[0,3] = 1.047416283515993
$inn[0,4] = 1.057296191440882
$inn[0,5] = 1.042688261945616
$ws.Range("I2:N2").Value = $inn

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.037051341506283
$bf[0,2] = 1.04304969165982
$bf[0,3] = 1.04548100243015
$bf[0,4] = 1.055446676686627
$ws.Range("B3:F3").Value = $bf

$inn = New-Object 'object[,]' 1,6
$inn[0,0] = 1.034314746747185
$inn[0,1] = 1.041805317070951
$inn[0,2] = 1.045637732967678
$inn[0,3] = 1.048062686062535
$inn[0,4] = 1.058002631086925
$inn[0,5] = 1.043284800117932
$ws.Range("I3:N3").Value = $inn

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.03766789880341
$bf[0,2] = 1.043510611168711
$bf[0,3] = 1.046022330412263
$bf[0,4] = 1.056025668133737
$ws.Range("B4:F4").Value = $bf

$inn = New-Object 'object[,]' 1,6
$inn[0,0] = 1.034397995007858
$inn[0,1] = 1.042190780797846
$inn[0,2] = 1.045975729512719
$inn[0,3] = 1.048481211112893
$inn[0,4] = 1.058460022975323
$inn[0,5] = 1.043670811247532
$ws.Range("I4:N4").Value = $inn

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.037927248415337
$bf[0,2] = 1.043704462010726
$bf[0,3] = 1.046250120438522
$bf[0,4] = 1.056269293787095
$ws.Range("B5:F5").Value = $bf

$inn = New-Object 'object[,]' 1,6
$inn[0,0] = 1.034432692873686
$inn[0,1] = 1.042352831283936
$inn[0,2] = 1.046117742792241
$inn[0,3] = 1.04865721991717
$inn[0,4] = 1.058652375747404
$inn[0,5] = 1.0438330918639
$ws.Range("I5:N5").Value = $inn

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.037970803056532
$bf[0,2] = 1.043737015082356
$bf[0,3] = 1.046288379993235
$bf[0,4] = 1.05631021232217
$ws.Range("B6:F6").Value = $bf

$inn = New-Object 'object[,]' 1,6
$inn[0,0] = 1.034438501206334
$inn[0,1] = 1.042380040331548
$inn[0,2] = 1.046141582693901
$inn[0,3] = 1.048686776103979
$inn[0,4] = 1.058684676408806
$inn[0,5] = 1.04386033955148
$ws.Range("I6:N6").Value = $inn

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.037671363661859
$bf[0,2] = 1.043513201099668
$bf[0,3] = 1.046025373307934
$bf[0,4] = 1.056028922618468
$ws.Range("B7:F7").Value = $bf

$inn = New-Object 'object[,]' 1,6
$inn[0,0] = 1.034398459820797
$inn[0,1] = 1.042192946119725
$inn[0,2] = 1.045977627418772
$inn[0,3] = 1.048483562713496
$inn[0,4] = 1.058462592947931
$inn[0,5] = 1.043672979644417
$ws.Range("I7:N7").Value = $inn

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.036421022845272
$bf[0,2] = 1.042578380732271
$bf[0,3] = 1.044927883381115
$bf[0,4] = 1.054855029109671
$ws.Range("B8:F8").Value = $bf

$inn = New-Object 'object[,]' 1,6
$inn[0,0] = 1.034228551428194
$inn[0,1] = 1.041410938877406
$inn[0,2] = 1.045291641398005
$inn[0,3] = 1.047634683841571
$inn[0,4] = 1.057534877594437
$inn[0,5] = 1.042889861862116
$ws.Range("I8:N8").Value = $inn

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.034222613228477
$bf[0,2] = 1.040933794985236
$bf[0,3] = 1.043000983826575
$bf[0,4] = 1.052793565302057
$ws.Range("B9:F9").Value = $bf

$inn = New-Object 'object[,]' 1,6
$inn[0,0] = 1.033919607766131
$inn[0,1] = 1.040033075456011
$inn[0,2] = 1.044080361645283
$inn[0,3] = 1.046140891364582
$inn[0,4] = 1.055902316221401
$inn[0,5] = 1.04151004171662
$ws.Range("I9:N9").Value = $inn

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.032760334678796
$bf[0,2] = 1.039839303811794
$bf[0,3] = 1.04172119287101
$bf[0,4] = 1.051424116922279
$ws.Range("B10:F10").Value = $bf

$inn = New-Object 'object[,]' 1,6
$inn[0,0] = 1.033707263593124
$inn[0,1] = 1.039114654746738
$inn[0,2] = 1.043271240856821
$inn[0,3] = 1.045146474067222
$inn[0,4] = 1.054815497215719
$inn[0,5] = 1.040590316744595
$ws.Range("I10:N10").Value = $inn

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.032127954948503
$bf[0,2] = 1.039365846973082
$bf[0,3] = 1.041168188740511
$bf[0,4] = 1.050832303730831
$ws.Range("B11:F11").Value = $bf

$inn = New-Object 'object[,]' 1,6
$inn[0,0] = 1.033613809806231
$inn[0,1] = 1.038717018795786
$inn[0,2] = 1.042920517459058
$inn[0,3] = 1.044716238109599
$inn[0,4] = 1.054345278566805
$inn[0,5] = 1.040192116104983
$ws.Range("I11:N11").Value = $inn

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.03189318169165
$bf[0,2] = 1.039190055447115
$bf[0,3] = 1.040962953242969
$bf[0,4] = 1.050612654900608
$ws.Range("B12:F12").Value = $bf

$inn = New-Object 'object[,]' 1,6
$inn[0,0] = 1.033578871005018
$inn[0,1] = 1.038569327125573
$inn[0,2] = 1.042790188923266
$inn[0,3] = 1.044556483598621
$inn[0,4] = 1.054170677047868
$inn[0,5] = 1.040044214695659
$ws.Range("I12:N12").Value = $inn

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.031943535871733
$bf[0,2] = 1.039227760069118
$bf[0,3] = 1.041006969016832
$bf[0,4] = 1.050659762287532
$ws.Range("B13:F13").Value = $bf

$inn = New-Object 'object[,]' 1,6
$inn[0,0] = 1.0335863757101
$inn[0,1] = 1.038601007109701
$inn[0,2] = 1.042818147273032
$inn[0,3] = 1.044590749007946
$inn[0,4] = 1.05420812699169
$inn[0,5] = 1.040075939668997
$ws.Range("I13:N13").Value = $inn

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.032108546048802
$bf[0,2] = 1.039351314524684
$bf[0,3] = 1.041151220334434
$bf[0,4] = 1.050814143871168
$ws.Range("B14:F14").Value = $bf

$inn = New-Object 'object[,]' 1,6
$inn[0,0] = 1.033610926359149
$inn[0,1] = 1.038704810386294
$inn[0,2] = 1.042909745561654
$inn[0,3] = 1.044703031637641
$inn[0,4] = 1.054330844747009
$inn[0,5] = 1.04017989035815
$ws.Range("I14:N14").Value = $inn

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.032210230294899
$bf[0,2] = 1.039427450013492
$bf[0,3] = 1.041240121547907
$bf[0,4] = 1.050909286954517
$ws.Range("B15:F15").Value = $bf

$inn = New-Object 'object[,]' 1,6
$inn[0,0] = 1.033626022906634
$inn[0,1] = 1.038768768106714
$inn[0,2] = 1.042966175136216
$inn[0,3] = 1.04477221990016
$inn[0,4] = 1.054406463015164
$inn[0,5] = 1.040243938905868
$ws.Range("I15:N15").Value = $inn

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.032802320511313
$bf[0,2] = 1.039870735513908
$bf[0,3] = 1.0417579183723
$bf[0,4] = 1.051463418339504
$ws.Range("B16:F16").Value = $bf

$inn = New-Object 'object[,]' 1,6
$inn[0,0] = 1.033713434102624
$inn[0,1] = 1.039141045601878
$inn[0,2] = 1.043294509542661
$inn[0,3] = 1.045175034956825
$inn[0,4] = 1.054846712225328
$inn[0,5] = 1.040616745077776
$ws.Range("I16:N16").Value = $inn

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.033173937063238
$bf[0,2] = 1.040148922404801
$bf[0,3] = 1.042083028613513
$bf[0,4] = 1.051811323892521
$ws.Range("B17:F17").Value = $bf

$inn = New-Object 'object[,]' 1,6
$inn[0,0] = 1.033767861668639
$inn[0,1] = 1.039374578586691
$inn[0,2] = 1.043500367127546
$inn[0,3] = 1.045427805634515
$inn[0,4] = 1.055122971928808
$inn[0,5] = 1.040850609706215
$ws.Range("I17:N17").Value = $inn

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.033390771492931
$bf[0,2] = 1.040311228879267
$bf[0,3] = 1.04227277114008
$bf[0,4] = 1.05201436380536
$ws.Range("B18:F18").Value = $bf

$inn = New-Object 'object[,]' 1,6
$inn[0,0] = 1.033799462871181
$inn[0,1] = 1.039510798818601
$inn[0,2] = 1.043620404717012
$inn[0,3] = 1.045575276501428
$inn[0,4] = 1.055284146224348
$inn[0,5] = 1.040987023386479
$ws.Range("I18:N18").Value = $inn

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.033464719421324
$bf[0,2] = 1.040366578738737
$bf[0,3] = 1.042337487257603
$bf[0,4] = 1.05208361419982
$ws.Range("B19:F19").Value = $bf

$inn = New-Object 'object[,]' 1,6
$inn[0,0] = 1.033810213372353
$inn[0,1] = 1.039557247109896
$inn[0,2] = 1.043661328349903
$inn[0,3] = 1.045625565987722
$inn[0,4] = 1.055339108674968
$inn[0,5] = 1.041033537639675
$ws.Range("I19:N19").Value = $inn

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.03313405817452
$bf[0,2] = 1.040119070942389
$bf[0,3] = 1.042048135846152
$bf[0,4] = 1.051773985272221
$ws.Range("B20:F20").Value = $bf

$inn = New-Object 'object[,]' 1,6
$inn[0,0] = 1.033762037149109
$inn[0,1] = 1.039349522245401
$inn[0,2] = 1.043478284236743
$inn[0,3] = 1.045400682206787
$inn[0,4] = 1.055093328072843
$inn[0,5] = 1.040825517782046
$ws.Range("I20:N20").Value = $inn

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.032059951323086
$bf[0,2] = 1.039314928832889
$bf[0,3] = 1.041108737077646
$bf[0,4] = 1.050768677445301
$ws.Range("B21:F21").Value = $bf

$inn = New-Object 'object[,]' 1,6
$inn[0,0] = 1.03360370303646
$inn[0,1] = 1.038674242678358
$inn[0,2] = 1.042882773636995
$inn[0,3] = 1.044669965693722
$inn[0,4] = 1.054294705797501
$inn[0,5] = 1.040149279240562
$ws.Range("I21:N21").Value = $inn

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.031385316760459
$bf[0,2] = 1.038809746781556
$bf[0,3] = 1.040519112217263
$bf[0,4] = 1.050137625030105
$ws.Range("B22:F22").Value = $bf

$inn = New-Object 'object[,]' 1,6
$inn[0,0] = 1.033502845364268
$inn[0,1] = 1.038249714585532
$inn[0,2] = 1.04250803918435
$inn[0,3] = 1.044210849862197
$inn[0,4] = 1.05379292012841
$inn[0,5] = 1.03972414826915
$ws.Range("I22:N22").Value = $inn

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.031742886689099
$bf[0,2] = 1.03907751358336
$bf[0,3] = 1.040831587039607
$bf[0,4] = 1.050472060204614
$ws.Range("B23:F23").Value = $bf

$inn = New-Object 'object[,]' 1,6
$inn[0,0] = 1.033556435600356
$inn[0,1] = 1.038474760148404
$inn[0,2] = 1.042706722320748
$inn[0,3] = 1.044454205625417
$inn[0,4] = 1.054058893537329
$inn[0,5] = 1.039949513422534
$ws.Range("I23:N23").Value = $inn

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.033152077490598
$bf[0,2] = 1.04013255939516
$bf[0,3] = 1.042063902042285
$bf[0,4] = 1.051790856640112
$ws.Range("B24:F24").Value = $bf

$inn = New-Object 'object[,]' 1,6
$inn[0,0] = 1.033764669448324
$inn[0,1] = 1.039360844114023
$inn[0,2] = 1.043488262654841
$inn[0,3] = 1.04541293801131
$inn[0,4] = 1.055106722742142
$inn[0,5] = 1.04083685572902
$ws.Range("I24:N24").Value = $inn

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.034790373143384
$bf[0,2] = 1.04135863182902
$bf[0,3] = 1.043498292586116
$bf[0,4] = 1.053325653403846
$ws.Range("B25:F25").Value = $bf

$inn = New-Object 'object[,]' 1,6
$inn[0,0] = 1.034000604475802
$inn[0,1] = 1.040389263660504
$inn[0,2] = 1.044393793774605
$inn[0,3] = 1.046526823251818
$inn[0,4] = 1.056324104332246
$inn[0,5] = 1.041866735749218
$ws.Range("I25:N25").Value = $inn

Write-Host "Applied vm_pu updates for rows 2-25"
